# Add files via upload
# Inserts a new worksheet ("Sheet2") after "Sheet1" holding a
# Reconciliation Deviation Scenarios explanation table, and leaves it as
# the active / selected sheet (matching the workbook's activeTab=1).

$wb = $excel.ActiveWorkbook
$sheet1 = $wb.Worksheets.Item(1)

$ws2 = $wb.Worksheets.Add($null, $sheet1)
$ws2.Name = "Sheet2"

# ---- Title (B3) ----
$ws2.Range("B3").Value = "📊 Reconciliation Deviation Scenarios – Detailed Business Explanation"
$ws2.Range("B3").Font.Bold = $true
$ws2.Range("B3").Font.Size = 18

# ---- Header row (row 5) ----
$ws2.Range("B5").Value = "Scenario"
$ws2.Range("C5").Value = "What the report shows"
$ws2.Range("D5").Value = "Why this happens (Detailed)"
$ws2.Range("E5").Value = "Environment scope"
$ws2.Range("F5").Value = "Should this be flagged as an exception?"
$ws2.Range("G5").Value = "Business decision needed"
$ws2.Range("B5:G5").Font.Bold = $true

# ---- Row 6: Matched in both (Backdated) ----
$ws2.Range("B6").Value = "Matched in both (Backdated)"
$ws2.Range("B6").Font.Bold = $true
$ws2.Range("C6").Value = "Record exists in Athena and MQS, but timestamp difference exceeds expected window"

$ws2.Range("D6").Value = "In test environments, the application team intentionally backdates quotes to simulate renewal journeys. This changes MQS/Athena timestamps but does not represent data loss."
$ws2.Range("D6").Characters(44, 30).Font.Bold = $true   # "intentionally backdates quotes"
$ws2.Range("D6").Characters(87, 16).Font.Bold = $true   # "renewal journeys"

$ws2.Range("E6").Value = "Test only"
$ws2.Range("E6").Font.Bold = $true

$ws2.Range("F6").Value = "⚠️ Optional – recommend reporting separately, not failing reconciliation"
$ws2.Range("F6").Characters(4, 8).Font.Italic = $true   # "Optional"

$ws2.Range("G6").Value = "Decide whether this should be an informational exception in test reports"
$ws2.Range("G6").Characters(34, 23).Font.Bold = $true   # "informational exception"

# ---- Row 7: MQS only ----
$ws2.Range("B7").Value = "MQS only"
$ws2.Range("B7").Font.Bold = $true
$ws2.Range("C7").Value = "Record exists in MQS but missing in Athena"

$ws2.Range("D7").Value = "Caused by Athena environment instability or ingestion issues, where events are not published or consumed correctly. This results in analytics data loss, not business transaction loss."
$ws2.Range("D7").Characters(11, 30).Font.Bold = $true    # "Athena environment instability"
$ws2.Range("D7").Characters(133, 19).Font.Bold = $true   # "analytics data loss"

$ws2.Range("E7").Value = "Test & non-prod"

$ws2.Range("F7").Value = "❌ Yes – real deviation"
$ws2.Range("F7").Characters(3, 20).Font.Bold = $true   # "Yes – real deviation"

$ws2.Range("G7").Value = "Requires Data Engineering / Athena environment remediation"

# ---- Row 8: CSV quoteId missing MQS prefix ----
$ws2.Range("B8").Value = "CSV quoteId missing MQS prefix"
$ws2.Range("B8").Font.Bold = $true

$ws2.Range("C8").Value = "QuoteId present but without MQS- prefix"
$c8 = $ws2.Range("C8").Characters(29, 4)   # "MQS-"
$c8.Font.Name = "Arial Unicode MS"
$c8.Font.Size = 10

$ws2.Range("D8").Value = "Used by teams to distinguish mock services from real services during testing. This affects reconciliation logic but not real business processing."
$ws2.Range("D8").Characters(18, 44).Font.Bold = $true   # "distinguish mock services from real services"

$ws2.Range("E8").Value = "Test only"
$ws2.Range("E8").Font.Bold = $true

$ws2.Range("F8").Value = "⚠️ No (exclude from prod-style KPI)"
$ws2.Range("F8").Characters(4, 2).Font.Italic = $true   # "No"

$ws2.Range("G8").Value = "Agree that this is a test-only data quality artefact"
$ws2.Range("G8").Characters(22, 31).Font.Bold = $true   # "test-only data quality artefact"

# ---- Row 9: CSV quoteId is NULL ----
$ws2.Range("B9").Value = "CSV quoteId is NULL"
$ws2.Range("B9").Font.Bold = $true
$ws2.Range("C9").Value = "CSV record has no quoteId"

$ws2.Range("D9").Value = "Certain products or test scenarios do not generate a quoteId, or source extract is incomplete. These rows cannot be reconciled technically."
$ws2.Range("D9").Characters(36, 25).Font.Bold = $true   # "do not generate a quoteId"

$ws2.Range("E9").Value = "Test & prod (product-dependent)"
$ws2.Range("F9").Value = "⚠️ Informational"
$ws2.Range("G9").Value = "Exclude from reconciliation denominator"

# ---- Row 10: Athena only ----
$ws2.Range("B10").Value = "Athena only"
$ws2.Range("B10").Font.Bold = $true
$ws2.Range("C10").Value = "Record exists in Athena but missing in MQS"
$ws2.Range("D10").Value = "MQS Product Adapter not triggered or still pointing to legacy DB during transition testing"
$ws2.Range("E10").Value = "Test / transition phase"
$ws2.Range("F10").Value = "⚠️ Conditional"
$ws2.Range("G10").Value = "MQS team confirmation required"

# ---- Row 11: Missing in both Athena & MQS ----
$ws2.Range("B11").Value = "Missing in both Athena & MQS"
$ws2.Range("B11").Font.Bold = $true
$ws2.Range("C11").Value = "CSV record not found anywhere"
$ws2.Range("D11").Value = "Source system did not publish event to EventHub, or event was dropped upstream"
$ws2.Range("E11").Value = "All environments"

$ws2.Range("F11").Value = "❌ Critical"
$ws2.Range("F11").Characters(3, 8).Font.Bold = $true   # "Critical"

$ws2.Range("G11").Value = "Immediate investigation required"

# ---- Selection / view state ----
$ws2.Range("B3:G11").Select()

